# PIR Sensor and Index.html
# Applies the "Questions" column + row 17 ("Boxing") additions to the
# Steps sheet, corrects the "5//5" / "6//6" typos on the Overview sheet,
# adds a hyperlink from the new Questions column back into the Steps
# sheet, and leaves the Steps sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: fix the date-look-alike text entries (re-enter with a
# leading apostrophe so Excel keeps them as literal text against the
# existing "m/d" number format already applied to these cells).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B5").Value = "'5/5"
$wsOverview.Range("B6").Value = "'6/6"

# ---------------------------------------------------------------------
# Steps sheet: add the new "Questions" column (D) and a new row (17)
# about boxing/packaging the project.
# ---------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("Steps")

$wsSteps.Columns.Item(4).ColumnWidth = 38.21875

$wsSteps.Range("D1").Value = "Questions"
$wsSteps.Range("D2").Value = "How to connect? Do I put both wires next to each other?"

$wsSteps.Range("D12").Value = "1. Put a image of the relay and explain the differnet pins and sockets. `n2. Then explain in a para how the relay works`n3. Then show the connection in 2 parts - part a from arduino to relay and part b - relay to pump"
$wsSteps.Range("D12").WrapText = $true
$wsSteps.Rows.Item(12).RowHeight = 100.8

$wsSteps.Range("D3").Value = "Give this as 1 code - Arduino + Sensor"
$wsSteps.Range("D7").Value = "Give this as 2nd code"
$wsSteps.Range("D11").Value = "Give this as 3rd code"

# D15 becomes a hyperlink back to the explanation in D12, displayed as
# "Answered here".
$wsSteps.Hyperlinks.Add($wsSteps.Range("D15"), "", "Steps!D12", "", "Answered here")

$wsSteps.Range("D16").Value = "Give as code 4"

# New row 17: "Boxing" step.
$wsSteps.Range("A17").Value = 7
$wsSteps.Range("B17").Value = "Boxing"
$wsSteps.Range("C17").Value = "Show how to package all this in one box"

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping to match the edited file: Steps
# becomes the active sheet (selected at A17); Overview's remembered
# selection moves to B7.
# ---------------------------------------------------------------------
$wsOverview.Range("B7").Select()
$wsSteps.Activate()
$wsSteps.Range("A17").Select()
